# Applied css to profile page
# Adds two new boolean tracking columns ("Leave_Applied" in N, "Leave_Approved" in O)
# to the student roster sheet, sets sample values for the two data rows, widens the
# new "Leave_Applied" column, and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells. Shared-string table order in the target workbook has
# "Leave_Approved" registered before "Leave_Applied", so write O1 first.
$ws.Range("O1").Value = "Leave_Approved"
$ws.Range("N1").Value = "Leave_Applied"

# Match the header look (bold / centered / wrapped) already used by the rest of row 1.
$ws.Range("M1").Copy() | Out-Null
$ws.Range("N1:O1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New per-row boolean values.
$ws.Range("N2").Value = $true
$ws.Range("O2").Value = $true

$ws.Range("N3").Value = $false
$ws.Range("O3").Value = $true

# Widen the new Leave_Applied column.
$ws.Range("N1").ColumnWidth = 15.14

# Move the viewport/selection the way the author left it.
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("O4").Select() | Out-Null
